# Auto-generated script applying the Mateus_Profits.xlsx market-data refresh
# (scheduled runner updated currentAveragePrice / Leve price / profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9776.786
$ws.Range("I64").Value = 3645.8333
$ws.Range("J64").Value = 14375
$ws.Range("K64").Value = 3645.8333
$ws.Range("L64").Value = 14375
$ws.Range("M64").Value = -3397.8333
$ws.Range("N64").Value = -14871
$ws.Range("H67").Value = 9776.786
$ws.Range("I67").Value = 3645.8333
$ws.Range("J67").Value = 14375
$ws.Range("K67").Value = 3645.8333
$ws.Range("L67").Value = 14375
$ws.Range("M67").Value = -2787.8333
$ws.Range("N67").Value = -16091
$ws.Range("H70").Value = 3219.6
$ws.Range("I70").Value = 2199
$ws.Range("J70").Value = 3900
$ws.Range("K70").Value = 6597
$ws.Range("L70").Value = 11700
$ws.Range("M70").Value = -6327
$ws.Range("N70").Value = -12240
$ws.Range("H73").Value = 3219.6
$ws.Range("I73").Value = 2199
$ws.Range("J73").Value = 3900
$ws.Range("K73").Value = 6597
$ws.Range("L73").Value = 11700
$ws.Range("M73").Value = -5661
$ws.Range("N73").Value = -13572
$ws.Range("H97").Value = 917.6429000000001
$ws.Range("J97").Value = 917.6429000000001
$ws.Range("L97").Value = 2752.9287
$ws.Range("N97").Value = -3744.9287
$ws.Range("H132").Value = 7145.3335
$ws.Range("I132").Value = 1154.7097
$ws.Range("K132").Value = 3464.1291
$ws.Range("M132").Value = -934.1291000000001
$ws.Range("H138").Value = 5407.5835
$ws.Range("I138").Value = 3046
$ws.Range("J138").Value = 5744.952
$ws.Range("K138").Value = 9138
$ws.Range("L138").Value = 17234.856
$ws.Range("M138").Value = -3998
$ws.Range("N138").Value = -27514.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14712646
$ws.Range("I61").Value = 15631561
$ws.Range("K61").Value = 15631561
$ws.Range("M61").Value = -15631349
$ws.Range("H74").Value = 2826.862
$ws.Range("I74").Value = 1769.05
$ws.Range("J74").Value = 5177.5557
$ws.Range("K74").Value = 1769.05
$ws.Range("L74").Value = 5177.5557
$ws.Range("M74").Value = -895.05
$ws.Range("N74").Value = -6925.5557
$ws.Range("H77").Value = 2826.862
$ws.Range("I77").Value = 1769.05
$ws.Range("J77").Value = 5177.5557
$ws.Range("K77").Value = 8845.25
$ws.Range("L77").Value = 25887.7785
$ws.Range("M77").Value = -4477.25
$ws.Range("N77").Value = -34623.7785
$ws.Range("H136").Value = 14712646
$ws.Range("I136").Value = 15631561
$ws.Range("K136").Value = 46894683
$ws.Range("M136").Value = -46892133

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2134.5293
$ws.Range("I22").Value = 2352.4666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 2352.4666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -2179.4666
$ws.Range("N22").Value = -846
$ws.Range("H60").Value = 149000
$ws.Range("J60").Value = 149000
$ws.Range("L60").Value = 149000
$ws.Range("N60").Value = -150198
$ws.Range("H134").Value = 4753.5127
$ws.Range("I134").Value = 4839.1313
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 14517.3939
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -11982.3939
$ws.Range("N134").Value = -9570
$ws.Range("H135").Value = 89994.5
$ws.Range("J135").Value = 89994.5
$ws.Range("L135").Value = 89994.5
$ws.Range("N135").Value = -100134.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 88191
$ws.Range("J52").Value = 88191
$ws.Range("L52").Value = 88191
$ws.Range("N52").Value = -88779
$ws.Range("H62").Value = 9499.375
$ws.Range("I62").Value = 9499.166999999999
$ws.Range("K62").Value = 9499.166999999999
$ws.Range("M62").Value = -8875.166999999999
$ws.Range("H65").Value = 9499.375
$ws.Range("I65").Value = 9499.166999999999
$ws.Range("K65").Value = 47495.835
$ws.Range("M65").Value = -44375.835
$ws.Range("H122").Value = 3664.6875
$ws.Range("I122").Value = 3853.5
$ws.Range("K122").Value = 11560.5
$ws.Range("M122").Value = -9110.5
$ws.Range("H132").Value = 2426.9092
$ws.Range("I132").Value = 2426.9092
$ws.Range("K132").Value = 7280.7276
$ws.Range("M132").Value = -4750.7276
$ws.Range("H133").Value = 44516.555
$ws.Range("J133").Value = 48807
$ws.Range("L133").Value = 48807
$ws.Range("N133").Value = -53867
$ws.Range("H134").Value = 4031.4595
$ws.Range("I134").Value = 3201.6562
$ws.Range("K134").Value = 9604.9686
$ws.Range("M134").Value = -7069.9686
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1173.25
$ws.Range("J5").Value = 1899.5
$ws.Range("L5").Value = 5698.5
$ws.Range("N5").Value = -5922.5
$ws.Range("H37").Value = 112733.8
$ws.Range("J37").Value = 112733.8
$ws.Range("L37").Value = 338201.4
$ws.Range("N37").Value = -338425.4
$ws.Range("H56").Value = 4999.8335
$ws.Range("I56").Value = 4999.8335
$ws.Range("K56").Value = 4999.8335
$ws.Range("M56").Value = -4469.8335
$ws.Range("H122").Value = 1546.375
$ws.Range("I122").Value = 796.6667
$ws.Range("K122").Value = 7170.0003
$ws.Range("M122").Value = -4720.0003
$ws.Range("H135").Value = 1173.25
$ws.Range("J135").Value = 1899.5
$ws.Range("L135").Value = 17095.5
$ws.Range("N135").Value = -22165.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 96504410
$ws.Range("I11").Value = 151504930
$ws.Range("K11").Value = 151504930
$ws.Range("M11").Value = -151504791
$ws.Range("H18").Value = 7667
$ws.Range("I18").Value = 4000.5
$ws.Range("J18").Value = 15000
$ws.Range("K18").Value = 4000.5
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = -3707.5
$ws.Range("N18").Value = -15586
$ws.Range("H80").Value = 2857.1765
$ws.Range("I80").Value = 3219.75
$ws.Range("K80").Value = 3219.75
$ws.Range("M80").Value = -2221.75
$ws.Range("H83").Value = 2857.1765
$ws.Range("I83").Value = 3219.75
$ws.Range("K83").Value = 16098.75
$ws.Range("M83").Value = -11106.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4518.154
$ws.Range("I40").Value = 4105.143
$ws.Range("K40").Value = 4105.143
$ws.Range("M40").Value = -3969.143
$ws.Range("H93").Value = 7449.7437
$ws.Range("I93").Value = 2157.037
$ws.Range("J93").Value = 19358.334
$ws.Range("K93").Value = 2157.037
$ws.Range("L93").Value = 19358.334
$ws.Range("M93").Value = -909.0369999999998
$ws.Range("N93").Value = -21854.334
$ws.Range("H136").Value = 6561.2
$ws.Range("I136").Value = 7201.5
$ws.Range("K136").Value = 21604.5
$ws.Range("M136").Value = -19054.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802
$ws.Range("H113").Value = 489.75
$ws.Range("I113").Value = 354.65
$ws.Range("J113").Value = 1165.25
$ws.Range("K113").Value = 1063.95
$ws.Range("L113").Value = 3495.75
$ws.Range("M113").Value = 1106.05
$ws.Range("N113").Value = -7835.75
$ws.Range("H141").Value = 74948.75
$ws.Range("J141").Value = 74948.75
$ws.Range("L141").Value = 74948.75
$ws.Range("N141").Value = -85308.75

